$d = $word.ActiveDocument

# Locate the paragraph that ends with "Many partial (some of the codes independently)"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Many partial \(some of the codes independently\)") {
        $target = $p
    }
}

$endRange = $target.Range
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Range.Text = "Legal aspects"
$newPara.Range.ListFormat.ListLevelNumber = 1
